$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 163; this shifts the existing rows 163-252
# down to 164-253 (preserving all of their data/formatting), matching the
# diff which shows every row from 164 to 253 now holding the data that
# used to belong to the row directly above it.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new weekly price record.
$ws.Range("A163").Value2 = 7
$ws.Range("B163").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C163").Value2 = "Ñuble"
$ws.Range("D163").Value2 = 44777
$ws.Range("E163").Value2 = 16
$ws.Range("F163").Value2 = 100112006
$ws.Range("G163").Value2 = "Repollo"
$ws.Range("H163").Value2 = "Crespo record"
$ws.Range("I163").Value2 = "Primera"
$ws.Range("J163").Value2 = 120
$ws.Range("K163").Value2 = 1200
$ws.Range("L163").Value2 = 1300
$ws.Range("M163").Value2 = 1250
$ws.Range("N163").Value2 = "$/unidad"
$ws.Range("O163").Value2 = "Provincia de Diguillín"
$ws.Range("P163").Value2 = 1250
$ws.Range("Q163").Value2 = 1
$ws.Range("R163").Value2 = "Hortaliza"
